$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (column G) for first data row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 03:03:19"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 03:03:15"
$wsZhCn.Range("K2").Value = "2016-08-23 03:03:39"

# de-de sheet: Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 03:03:19"
$wsDeDe.Range("K2").Value = "2016-08-23 03:03:46"
